$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15
$ws1.Range("F3").Value = 971
$ws1.Range("F6").Value = 1088
$ws1.Range("F7").Value = 872
$ws1.Range("F8").Value = 263
$ws1.Range("F13").Value = 588
$ws1.Range("F15").Value = 1352
$ws1.Range("F18").Value = 1224
$ws1.Range("F20").Value = 1487
$ws1.Range("F21").Value = 727
$ws1.Range("F23").Value = 1293
$ws1.Range("F25").Value = 1042
$ws1.Range("F27").Value = 3211
$ws1.Range("F28").Value = 626
$ws1.Range("F29").Value = 539

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 161
$ws2.Range("F8").Value = 29

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 757

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15
$ws4.Range("F3").Value = 757
$ws4.Range("F6").Value = 971
$ws4.Range("F10").Value = 1088
$ws4.Range("F11").Value = 872
$ws4.Range("F12").Value = 263
$ws4.Range("F17").Value = 161
$ws4.Range("F18").Value = 29
$ws4.Range("F23").Value = 588
$ws4.Range("F25").Value = 1352
$ws4.Range("F28").Value = 1224
$ws4.Range("F30").Value = 1487
$ws4.Range("F31").Value = 727
$ws4.Range("F33").Value = 1293
$ws4.Range("F37").Value = 1042
$ws4.Range("F39").Value = 3211
$ws4.Range("F40").Value = 626
$ws4.Range("F41").Value = 539
